$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 162, shifting existing rows 162-238 down
# to 163-239 (dimension grows from A1:R238 to A1:R239).
$ws.Rows(162).Insert()

# Populate the new row 162 with the new data point (same Mercado/Region/
# Categoria/Calidad/Unidad/Origen/Kg-unidades/Clasificacion as the
# surrounding rows, with a new Fecha/Volumen/Precios).
$ws.Range("A162").Value = 3
$ws.Range("B162").Value = "Femacal de La Calera"
$ws.Range("C162").Value = "Coquimbo"
$ws.Range("D162").Value = 44510
$ws.Range("E162").Value = 5
$ws.Range("F162").Value = 100112040
$ws.Range("G162").Value = "Cilantro"
$ws.Range("H162").Value = "Sin especificar"
$ws.Range("I162").Value = "Primera"
$ws.Range("J162").Value = 240
$ws.Range("K162").Value = 2000
$ws.Range("L162").Value = 2500
$ws.Range("M162").Value = 2250
$ws.Range("N162").Value = "$/docena de atados (3 kilos)"
$ws.Range("O162").Value = "Provincia de Quillota"
$ws.Range("P162").Value = 750
$ws.Range("Q162").Value = 3
$ws.Range("R162").Value = "Hortaliza"
